# Adiciona o campo "Sujeito a ISS?" (nova coluna U), deslocando as colunas
# "Outros Impostos" (U -> V) e "Info Adicionais" (V -> W) para a direita.
# Tambem limpa alguns campos de teste ("_BLOCO") que haviam ficado no JSON
# da coluna C e marca "SIM" nas linhas cujo item possua um valor em
# ISSQN_BLOCO.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insere uma nova coluna antes da coluna U (desloca U->V, V->W, ...)
$ws.Columns("U:U").Insert()

# 2) Cabecalho da nova coluna
$ws.Range("U1").Value = "Sujeito a ISS?"

# 3) Preenche a nova coluna "Sujeito a ISS?" com "NAO" para todas as linhas
#    de dados (linha 1 e o cabecalho).
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 21).Value = "NAO"
}

# 4) Linhas cujo item possui valor em ISSQN_BLOCO passam a ser "SIM"
$ws.Cells.Item(20, 21).Value = "SIM"
$ws.Cells.Item(30, 21).Value = "SIM"

# 5) Remove campos de teste/residuais do JSON (coluna C) e ajusta
#    ISSQN_BLOCO para as linhas 20 e 30.

# --- Linha 2 ---
$v = $ws.Range("C2").Value2
$v = $v.Replace("`"INFO_ADICIONAL`": null,`n      `"II_BLOCO`": `"II`"", "`"INFO_ADICIONAL`": null")
$v = $v.Replace("`"INFO_ADICIONAL`": null,`n      `"IOF_BLOCO`": `"IOF`"", "`"INFO_ADICIONAL`": null")
$ws.Range("C2").Value = $v

# --- Linha 3 (mesmo JSON da linha 2) ---
$v = $ws.Range("C3").Value2
$v = $v.Replace("`"INFO_ADICIONAL`": null,`n      `"II_BLOCO`": `"II`"", "`"INFO_ADICIONAL`": null")
$v = $v.Replace("`"INFO_ADICIONAL`": null,`n      `"IOF_BLOCO`": `"IOF`"", "`"INFO_ADICIONAL`": null")
$ws.Range("C3").Value = $v

# --- Linha 4 ---
$v = $ws.Range("C4").Value2
$v = $v.Replace("`"IPI_BLOCO`": `"teste`",", "`"IPI_BLOCO`": null,")
$v = $v.Replace("`"PIS_BLOCO`": `"TESTEISSQN`",`n      `"ISSQN_BLOCO`": `"TESTEISSQN`",", "`"PIS_BLOCO`": `"PISNT`",`n      `"ISSQN_BLOCO`": null,")
$v = $v.Replace("`"OUTRO_IMPOSTO`": `"IOF`",`n      `"ICMS_UFDEST_BLOCO`": `"TESTEISSQN`",", "`"OUTRO_IMPOSTO`": null,`n      `"ICMS_UFDEST_BLOCO`": null,")
$v = $v.Replace("`"INFO_ADICIONAL`": null,`n      `"CIDE_BLOCO`": `"CIDE`"", "`"INFO_ADICIONAL`": null")
$ws.Range("C4").Value = $v

# --- Linha 5 ---
$v = $ws.Range("C5").Value2
$v = $v.Replace("`"INFO_ADICIONAL`": null,`n      `"II_BLOCO`": `"II`",`n      `"IOF_BLOCO`": `"IOF`",`n      `"CIDE_BLOCO`": `"CIDE`"", "`"INFO_ADICIONAL`": null")
$ws.Range("C5").Value = $v

# --- Linha 6 ---
$v = $ws.Range("C6").Value2
$v = $v.Replace("`"ICMS_UFDEST_BLOCO`": `"ICMSDUFDEST`",", "`"ICMS_UFDEST_BLOCO`": null,")
$v = $v.Replace("`"INFO_ADICIONAL`": null,`n      `"PMPF_BLOCO`": `"PMPF`"", "`"INFO_ADICIONAL`": null")
$ws.Range("C6").Value = $v

# --- Linha 8 ---
$v = $ws.Range("C8").Value2
$v = $v.Replace("`"ICMS_UFDEST_BLOCO`": `"ICMSDUFDEST`",", "`"ICMS_UFDEST_BLOCO`": null,")
$ws.Range("C8").Value = $v

# --- Linha 20 ---
$v = $ws.Range("C20").Value2
$v = $v.Replace("`"PIS_BLOCO`": `"PISNT`",`n      `"ISSQN_BLOCO`": null,", "`"PIS_BLOCO`": `"PISNT`",`n      `"ISSQN_BLOCO`": `"valeu`",")
$ws.Range("C20").Value = $v

# --- Linha 30 ---
$v = $ws.Range("C30").Value2
$v = $v.Replace("`"PIS_BLOCO`": `"PISNT`",`n      `"ISSQN_BLOCO`": null,", "`"PIS_BLOCO`": `"PISNT`",`n      `"ISSQN_BLOCO`": `"valeu`",")
$ws.Range("C30").Value = $v

Write-Output "done"
